$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3647.9656
$ws.Range("I40").Value = 1997.25
$ws.Range("J40").Value = 3912.08
$ws.Range("K40").Value = 1997.25
$ws.Range("L40").Value = 3912.08
$ws.Range("M40").Value = -1822.25
$ws.Range("N40").Value = -4262.08
$ws.Range("H43").Value = 6246.5
$ws.Range("I43").Value = 4997
$ws.Range("K43").Value = 4997
$ws.Range("M43").Value = -4928
$ws.Range("H48").Value = 11955
$ws.Range("I48").Value = 9500
$ws.Range("J48").Value = 13591.667
$ws.Range("K48").Value = 28500
$ws.Range("L48").Value = 40775.001
$ws.Range("M48").Value = -28208
$ws.Range("N48").Value = -41359.001
$ws.Range("H56").Value = 11955
$ws.Range("I56").Value = 9500
$ws.Range("J56").Value = 13591.667
$ws.Range("K56").Value = 28500
$ws.Range("L56").Value = 40775.001
$ws.Range("M56").Value = -27966
$ws.Range("N56").Value = -41843.001
$ws.Range("H131").Value = 2880.5881
$ws.Range("I131").Value = 796.8889
$ws.Range("K131").Value = 2390.6667
$ws.Range("M131").Value = 2649.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25394716
$ws.Range("I32").Value = 31084066
$ws.Range("K32").Value = 31084066
$ws.Range("M32").Value = -31083779
$ws.Range("H61").Value = 3556.5898
$ws.Range("I61").Value = 3504.2307
$ws.Range("K61").Value = 3504.2307
$ws.Range("M61").Value = -3292.2307
$ws.Range("H124").Value = 38500
$ws.Range("J124").Value = 38500
$ws.Range("L124").Value = 38500
$ws.Range("N124").Value = -48320
$ws.Range("H132").Value = 2746.1455
$ws.Range("I132").Value = 2465.587
$ws.Range("J132").Value = 4180.1113
$ws.Range("K132").Value = 7396.761
$ws.Range("L132").Value = 12540.3339
$ws.Range("M132").Value = -4866.761
$ws.Range("N132").Value = -17600.3339
$ws.Range("H136").Value = 3556.5898
$ws.Range("I136").Value = 3504.2307
$ws.Range("K136").Value = 10512.6921
$ws.Range("M136").Value = -7962.6921
$ws.Range("H137").Value = 97599.39999999999
$ws.Range("I137").Value = 43999.168
$ws.Range("J137").Value = 177999.75
$ws.Range("K137").Value = 43999.168
$ws.Range("L137").Value = 177999.75
$ws.Range("M137").Value = -38899.168
$ws.Range("N137").Value = -188199.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3084.3333
$ws.Range("I86").Value = 2925
$ws.Range("K86").Value = 2925
$ws.Range("M86").Value = -1802
$ws.Range("H89").Value = 3084.3333
$ws.Range("I89").Value = 2925
$ws.Range("K89").Value = 14625
$ws.Range("M89").Value = -9009
$ws.Range("H94").Value = 1513.7646
$ws.Range("I94").Value = 1435
$ws.Range("K94").Value = 1435
$ws.Range("M94").Value = -984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 58408.375
$ws.Range("I69").Value = 42999.5
$ws.Range("K69").Value = 42999.5
$ws.Range("M69").Value = -42250.5
$ws.Range("H72").Value = 58408.375
$ws.Range("I72").Value = 42999.5
$ws.Range("K72").Value = 128998.5
$ws.Range("M72").Value = -125254.5
$ws.Range("H132").Value = 3868.5676
$ws.Range("I132").Value = 3818.5
$ws.Range("K132").Value = 11455.5
$ws.Range("M132").Value = -8925.5
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060
$ws.Range("H134").Value = 2504.9333
$ws.Range("I134").Value = 2375.9092
$ws.Range("J134").Value = 2859.75
$ws.Range("K134").Value = 7127.7276
$ws.Range("L134").Value = 8579.25
$ws.Range("M134").Value = -4592.7276
$ws.Range("N134").Value = -13649.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1136.1875
$ws.Range("I5").Value = 793.7692
$ws.Range("K5").Value = 2381.3076
$ws.Range("M5").Value = -2269.3076
$ws.Range("H18").Value = 5174
$ws.Range("I18").Value = 4073
$ws.Range("K18").Value = 12219
$ws.Range("M18").Value = -12050
$ws.Range("H132").Value = 734.5714
$ws.Range("J132").Value = 1495
$ws.Range("L132").Value = 13455
$ws.Range("N132").Value = -18515
$ws.Range("H133").Value = 4951.9
$ws.Range("J133").Value = 5000
$ws.Range("L133").Value = 15000
$ws.Range("N133").Value = -25120
$ws.Range("H134").Value = 4515.6665
$ws.Range("I134").Value = 3520.3572
$ws.Range("K134").Value = 10561.0716
$ws.Range("M134").Value = -5491.071599999999
$ws.Range("H135").Value = 1136.1875
$ws.Range("I135").Value = 793.7692
$ws.Range("K135").Value = 7143.922799999999
$ws.Range("M135").Value = -4608.922799999999
$ws.Range("H136").Value = 2357.8
$ws.Range("I136").Value = 2357.8
$ws.Range("K136").Value = 7073.400000000001
$ws.Range("M136").Value = -1973.400000000001
$ws.Range("H137").Value = 23795.4
$ws.Range("I137").Value = 1326.6666
$ws.Range("K137").Value = 3979.9998
$ws.Range("M137").Value = 1120.0002
$ws.Range("H139").Value = 3653.6
$ws.Range("I139").Value = 2933
$ws.Range("K139").Value = 8799
$ws.Range("M139").Value = -3659

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 4792
$ws.Range("I2").Value = 32.125
$ws.Range("J2").Value = 20023.6
$ws.Range("K2").Value = 32.125
$ws.Range("L2").Value = 20023.6
$ws.Range("M2").Value = 80.875
$ws.Range("N2").Value = -20249.6
$ws.Range("H24").Value = 58843950
$ws.Range("I24").Value = 125011560
$ws.Range("J24").Value = 28295.777
$ws.Range("K24").Value = 125011560
$ws.Range("L24").Value = 28295.777
$ws.Range("M24").Value = -125011387
$ws.Range("N24").Value = -28641.777
$ws.Range("H44").Value = 50496.668
$ws.Range("I44").Value = 99990
$ws.Range("J44").Value = 25750
$ws.Range("K44").Value = 99990
$ws.Range("L44").Value = 25750
$ws.Range("M44").Value = -99394
$ws.Range("N44").Value = -26942
$ws.Range("H62").Value = 50000
$ws.Range("I62").Value = 44000
$ws.Range("J62").Value = 56000
$ws.Range("K62").Value = 44000
$ws.Range("L62").Value = 56000
$ws.Range("M62").Value = -43314
$ws.Range("N62").Value = -57372
$ws.Range("H65").Value = 50000
$ws.Range("I65").Value = 44000
$ws.Range("J65").Value = 56000
$ws.Range("K65").Value = 132000
$ws.Range("L65").Value = 168000
$ws.Range("M65").Value = -128568
$ws.Range("N65").Value = -174864
$ws.Range("H70").Value = 13296.892
$ws.Range("I70").Value = 82871.5
$ws.Range("K70").Value = 82871.5
$ws.Range("M70").Value = -82601.5
$ws.Range("H73").Value = 13296.892
$ws.Range("I73").Value = 82871.5
$ws.Range("K73").Value = 82871.5
$ws.Range("M73").Value = -81935.5
$ws.Range("H113").Value = 16186.571
$ws.Range("I113").Value = 2179
$ws.Range("K113").Value = 2179
$ws.Range("M113").Value = -9
$ws.Range("H132").Value = 4220.2085
$ws.Range("I132").Value = 4309.7896
$ws.Range("J132").Value = 3879.8
$ws.Range("K132").Value = 12929.3688
$ws.Range("L132").Value = 11639.4
$ws.Range("M132").Value = -10399.3688
$ws.Range("N132").Value = -16699.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3615.5
$ws.Range("I7").Value = 3489.1428
$ws.Range("J7").Value = 4500
$ws.Range("K7").Value = 3489.1428
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = -3377.1428
$ws.Range("N7").Value = -4724
$ws.Range("H122").Value = 44601.2
$ws.Range("I122").Value = 44001.5
$ws.Range("K122").Value = 132004.5
$ws.Range("M122").Value = -129554.5
$ws.Range("H126").Value = 3615.5
$ws.Range("I126").Value = 3489.1428
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 10467.4284
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -7997.428400000001
$ws.Range("N126").Value = -18440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1185.625
$ws.Range("I136").Value = 973.0333000000001
$ws.Range("K136").Value = 2919.0999
$ws.Range("M136").Value = -369.0999000000002
